$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Denmark Division 2")

function Swap-Rows($row1, $row2) {
    # Columns B (2) through AD (30) hold the data that should be swapped;
    # column A (1) is the row's fixed sequence number and must stay put.
    $firstCol = 2
    $lastCol = 30

    $rng1 = $ws.Range($ws.Cells.Item($row1, $firstCol), $ws.Cells.Item($row1, $lastCol))
    $rng2 = $ws.Range($ws.Cells.Item($row2, $firstCol), $ws.Cells.Item($row2, $lastCol))

    $vals1 = $rng1.Value()
    $vals2 = $rng2.Value()

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}

Swap-Rows 66 67
Swap-Rows 115 116
Swap-Rows 173 176
Swap-Rows 193 194
